# Updates cryptos list: refresh D (Price) and E (Volume(1h)) columns
# for rows 2-51 with newly scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column holds plain-text strings (e.g. "71.907.94", "11.70") in the
# source file; force text format first so Excel doesn't auto-coerce
# number-looking values (dropping trailing zeros, etc.) on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "71.907.94"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "3.717.30"
$ws.Range("E3").Value = "  +8.19%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "589.39"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "180.81"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("D7").Value = "3.703.97"
$ws.Range("E7").Value = "  +7.99%  "
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  +4.49%  "
$ws.Range("D12").Value = "50.05"
$ws.Range("E12").Value = "  +3.09%  "
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "4.315.72"
$ws.Range("E14").Value = "  +8.22%  "
$ws.Range("D15").Value = "683.54"
$ws.Range("E15").Value = "  -0.50%  "
$ws.Range("D16").Value = "9.06"
$ws.Range("E16").Value = "  +4.27%  "
$ws.Range("D17").Value = "3.722.03"
$ws.Range("E17").Value = "  +8.57%  "
$ws.Range("D18").Value = "71.947.63"
$ws.Range("E18").Value = "  +3.52%  "
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "18.19"
$ws.Range("E20").Value = "  +1.77%  "
$ws.Range("D21").Value = "11.70"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("D22").Value = "6.44"
$ws.Range("E22").Value = "  +19.33%  "
$ws.Range("D23").Value = "0.948"
$ws.Range("E23").Value = "  +3.55%  "
$ws.Range("D24").Value = "17.83"
$ws.Range("E24").Value = "  +4.63%  "
$ws.Range("D25").Value = "103.88"
$ws.Range("E25").Value = "  +2.62%  "
$ws.Range("E26").Value = "  +3.67%  "
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("D28").Value = "10.36"
$ws.Range("E28").Value = "  +6.57%  "
$ws.Range("D29").Value = "35.66"
$ws.Range("E29").Value = "  +5.77%  "
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("D31").Value = "7.38"
$ws.Range("E31").Value = "  +6.67%  "
$ws.Range("D32").Value = "4.23"
$ws.Range("E32").Value = "  +12.48%  "
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("E34").Value = "  +3.72%  "
$ws.Range("D35").Value = "562.95"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "59.69"
$ws.Range("E36").Value = "  +2.44%  "
$ws.Range("D37").Value = "3.777.85"
$ws.Range("E37").Value = "  +3.84%  "
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "0.144"
$ws.Range("E39").Value = "  +2.72%  "
$ws.Range("D40").Value = "0.0₃0781"
$ws.Range("E40").Value = "  +5.47%  "
$ws.Range("D41").Value = "35.87"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("D42").Value = "3.47"
$ws.Range("E42").Value = "  +5.38%  "
$ws.Range("E43").Value = "  +9.14%  "
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("E46").Value = "  +8.44%  "
$ws.Range("D47").Value = "3.39"
$ws.Range("E47").Value = "  +1.40%  "
$ws.Range("E48").Value = "  +3.49%  "
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "135.84"
$ws.Range("E51").Value = "  +3.55%  "
